# "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
#
# The previous "Estado de Cuenta" listed two workers (EDER DAVID GONZALEZ
# PALMERA for period 2409, and DANIEL ROYERO MARIN repeated across six
# periods 2507-2502). This update refreshes the account-statement database
# down to a single worker/period entry (part 1 of the refreshed data) and
# updates the summary totals to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the six rows (17-22) belonging to DANIEL ROYERO MARIN / periods
# 2507-2502. Excel shifts everything below up, so the signature block
# (previously rows 27-28) becomes rows 21-22, and the single remaining
# data row (EDER DAVID GONZALEZ PALMERA, period 2409) stays on row 16.
$ws.Rows("17:22").Delete() | Out-Null

# Update the summary header to reflect the now-single worker/period data:
# total "VALOR MORA" (E11), "Cant. Trabajadores" (C13) and "Cant. Periodos"
# (F13).
$ws.Range("E11").Value = 1733
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1
